$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels and Devices")

# Update the "Panel Accessories Devices" / "Panel Accessories Label" values for row 8
# (PCS800 / PCS800-1 removed from the comma separated lists)
$ws.Range("K8").Value = "FB800,IOB800,POS800-S,POS800-M"
$ws.Range("L8").Value = "FB800-1,IOB800-1,POS800-S-1,POS800-M-1"

# Update the expected 24V PSU load values for row 8
$ws.Range("J8").Value = 0.397
$ws.Range("N8").Value = 1.101
$ws.Range("O8").Value = 0.329

# Update the view state: scroll so column B is the left-most visible column,
# and select cell O12
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("O12").Select()
